$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format styling (style s="4") from the last existing data row (222)
# into the new rows' H/I columns before writing data, so the new cells reuse the
# existing shared cell style instead of creating a new one.
$ws.Range("H222:I222").Copy() | Out-Null
$ws.Range("H223:I226").PasteSpecial(-4122) | Out-Null

# Row 223: LeetCode 3190
$ws.Range("A223").Value = 3190
$ws.Range("B223").Value = "Find Minimum Operations to Make All Elements Divisible by Three"
$ws.Range("C223").Value = "#array"
$ws.Range("D223").Value = "easy"
$ws.Range("E223").Value = 1
$ws.Range("F223").Value = 0
$ws.Range("G223").Value = 2
$ws.Range("H223").Value = 45983
$ws.Range("I223").Value = 45983

# Row 224: LeetCode 1262
$ws.Range("A224").Value = 1262
$ws.Range("B224").Value = "Greatest Sum Divisible by Three"
$ws.Range("C224").Value = "#array #greedy #dynamic-programming #sorting "
$ws.Range("D224").Value = "medium"
$ws.Range("E224").Value = 0
$ws.Range("F224").Value = 1
$ws.Range("G224").Value = 55
$ws.Range("H224").Value = 45984
$ws.Range("I224").Value = 45984

# Row 225: LeetCode 1018
$ws.Range("A225").Value = 1018
$ws.Range("B225").Value = "Binary Prefix Divisible By 5"
$ws.Range("C225").Value = "#math "
$ws.Range("D225").Value = "easy"
$ws.Range("E225").Value = 1
$ws.Range("F225").Value = 0
$ws.Range("G225").Value = 5
$ws.Range("H225").Value = 45985
$ws.Range("I225").Value = 45985

# Row 226: LeetCode 364
$ws.Range("A226").Value = 364
$ws.Range("B226").Value = "Nested List Weight Sum II"
$ws.Range("C226").Value = "#bfs #dfs #deque #queue "
$ws.Range("D226").Value = "medium"
$ws.Range("E226").Value = 1
$ws.Range("F226").Value = 0
$ws.Range("G226").Value = 24
$ws.Range("H226").Value = 45986
$ws.Range("I226").Value = 45986

# Row heights: rows 223/224 wrap to two lines (ht=68), rows 225/226 to one line (ht=34)
$ws.Rows.Item(223).RowHeight = 68
$ws.Rows.Item(224).RowHeight = 68
$ws.Rows.Item(225).RowHeight = 34
$ws.Rows.Item(226).RowHeight = 34

# Update the view to match the end-of-sheet scroll/selection position recorded in the diff
$null = $ws.Range("J226").Select()
$excel.ActiveWindow.ScrollRow = 222
